# Auto-generated edit script: updates market-data-driven profit columns
# (currentAveragePrice*, LevePrice*, LeveProfit*) on several leve rows across
# all 8 job sheets, per the scheduled-runner refresh commit.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 69626.336
$ws.Range("J64").Value = 3299.5
$ws.Range("L64").Value = 3299.5
$ws.Range("N64").Value = -3795.5
$ws.Range("H67").Value = 69626.336
$ws.Range("J67").Value = 3299.5
$ws.Range("L67").Value = 3299.5
$ws.Range("N67").Value = -5015.5
$ws.Range("H70").Value = 1468.7778
$ws.Range("I70").Value = 1517.5883
$ws.Range("J70").Value = 1385.8
$ws.Range("K70").Value = 4552.7649
$ws.Range("L70").Value = 4157.4
$ws.Range("M70").Value = -4282.7649
$ws.Range("N70").Value = -4697.4
$ws.Range("H73").Value = 1468.7778
$ws.Range("I73").Value = 1517.5883
$ws.Range("J73").Value = 1385.8
$ws.Range("K73").Value = 4552.7649
$ws.Range("L73").Value = 4157.4
$ws.Range("M73").Value = -3616.7649
$ws.Range("N73").Value = -6029.4
$ws.Range("H108").Value = 39786.332
$ws.Range("J108").Value = 39786.332
$ws.Range("L108").Value = 39786.332
$ws.Range("N108").Value = -47466.332
$ws.Range("H138").Value = 3045.5217
$ws.Range("I138").Value = 1721
$ws.Range("J138").Value = 3595.7078
$ws.Range("K138").Value = 5163
$ws.Range("L138").Value = 10787.1234
$ws.Range("M138").Value = -23
$ws.Range("N138").Value = -21067.1234

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 3462.8333
$ws.Range("I26").Value = 2694.25
$ws.Range("K26").Value = 2694.25
$ws.Range("M26").Value = -2364.25
$ws.Range("H32").Value = 32127.812
$ws.Range("I32").Value = 11618.328
$ws.Range("K32").Value = 11618.328
$ws.Range("M32").Value = -11331.328
$ws.Range("H44").Value = 5988
$ws.Range("J44").Value = 5988
$ws.Range("L44").Value = 5988
$ws.Range("N44").Value = -6964
$ws.Range("H46").Value = 3604.4
$ws.Range("J46").Value = 3255.5
$ws.Range("L46").Value = 3255.5
$ws.Range("N46").Value = -3893.5
$ws.Range("H55").Value = 14490
$ws.Range("J55").Value = 15588
$ws.Range("L55").Value = 15588
$ws.Range("N55").Value = -16218
$ws.Range("H63").Value = 2743.4285
$ws.Range("J63").Value = 3550
$ws.Range("L63").Value = 3550
$ws.Range("N63").Value = -4922
$ws.Range("H64").Value = 43980.5
$ws.Range("J64").Value = 43980.5
$ws.Range("L64").Value = 43980.5
$ws.Range("N64").Value = -44476.5
$ws.Range("H66").Value = 2743.4285
$ws.Range("J66").Value = 3550
$ws.Range("L66").Value = 17750
$ws.Range("N66").Value = -24614
$ws.Range("H67").Value = 43980.5
$ws.Range("J67").Value = 43980.5
$ws.Range("L67").Value = 43980.5
$ws.Range("N67").Value = -45696.5
$ws.Range("H114").Value = 21333
$ws.Range("J114").Value = 21333
$ws.Range("L114").Value = 21333
$ws.Range("N114").Value = -30011
$ws.Range("H119").Value = 32000
$ws.Range("J119").Value = 32000
$ws.Range("L119").Value = 32000
$ws.Range("N119").Value = -41676
$ws.Range("H122").Value = 3454.3333
$ws.Range("I122").Value = 2527.4285
$ws.Range("J122").Value = 4752
$ws.Range("K122").Value = 7582.2855
$ws.Range("L122").Value = 14256
$ws.Range("M122").Value = -5132.2855
$ws.Range("N122").Value = -19156

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H42").Value = 334400.2
$ws.Range("J42").Value = 334400.2
$ws.Range("L42").Value = 334400.2
$ws.Range("N42").Value = -335056.2
$ws.Range("H62").Value = 45495
$ws.Range("J62").Value = 45495
$ws.Range("L62").Value = 45495
$ws.Range("N62").Value = -46867
$ws.Range("H65").Value = 45495
$ws.Range("J65").Value = 45495
$ws.Range("L65").Value = 136485
$ws.Range("N65").Value = -143349
$ws.Range("H105").Value = 223954.67
$ws.Range("I105").Value = 251895
$ws.Range("K105").Value = 251895
$ws.Range("M105").Value = -250148

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 733.75
$ws.Range("I22").Value = 206
$ws.Range("J22").Value = 1050.4
$ws.Range("K22").Value = 206
$ws.Range("L22").Value = 1050.4
$ws.Range("M22").Value = 144
$ws.Range("N22").Value = -1750.4
$ws.Range("H51").Value = 7897.2
$ws.Range("J51").Value = 7897.2
$ws.Range("L51").Value = 7897.2
$ws.Range("N51").Value = -9369.200000000001
$ws.Range("H58").Value = 1402.1389
$ws.Range("I58").Value = 1221.3438
$ws.Range("J58").Value = 2848.5
$ws.Range("K58").Value = 1221.3438
$ws.Range("L58").Value = 2848.5
$ws.Range("M58").Value = -1018.3438
$ws.Range("N58").Value = -3254.5
$ws.Range("H60").Value = 14528.6
$ws.Range("J60").Value = 14528.6
$ws.Range("L60").Value = 14528.6
$ws.Range("N60").Value = -15550.6
$ws.Range("H61").Value = 7897.2
$ws.Range("J61").Value = 7897.2
$ws.Range("L61").Value = 7897.2
$ws.Range("N61").Value = -8593.200000000001
$ws.Range("H95").Value = 20000
$ws.Range("J95").Value = 20000
$ws.Range("L95").Value = 20000
$ws.Range("N95").Value = -25492
$ws.Range("H107").Value = 840.125
$ws.Range("I107").Value = 803.2308
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 803.2308
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = 1116.7692
$ws.Range("N107").Value = -4840
$ws.Range("H133").Value = 25000
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H134").Value = 1242.3334
$ws.Range("I134").Value = 724.8
$ws.Range("J134").Value = 2536.1667
$ws.Range("K134").Value = 2174.4
$ws.Range("L134").Value = 7608.500100000001
$ws.Range("M134").Value = 360.6000000000004
$ws.Range("N134").Value = -12678.5001
$ws.Range("H136").Value = 1402.1389
$ws.Range("I136").Value = 1221.3438
$ws.Range("J136").Value = 2848.5
$ws.Range("K136").Value = 3664.0314
$ws.Range("L136").Value = 8545.5
$ws.Range("M136").Value = -1114.0314
$ws.Range("N136").Value = -13645.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1316.3636
$ws.Range("J34").Value = 1439
$ws.Range("L34").Value = 4317
$ws.Range("N34").Value = -4485
$ws.Range("H39").Value = 6276
$ws.Range("J39").Value = 8201.333000000001
$ws.Range("L39").Value = 24603.999
$ws.Range("N39").Value = -25191.999
$ws.Range("H55").Value = 9330
$ws.Range("J55").Value = 3205.5557
$ws.Range("L55").Value = 9616.667099999999
$ws.Range("N55").Value = -9970.667099999999
$ws.Range("H129").Value = 172456.4
$ws.Range("I129").Value = 6786.222
$ws.Range("J129").Value = 243457.9
$ws.Range("K129").Value = 20358.666
$ws.Range("L129").Value = 730373.7
$ws.Range("M129").Value = -15358.666
$ws.Range("N129").Value = -740373.7
$ws.Range("H131").Value = 803.57574
$ws.Range("I131").Value = 520.64703
$ws.Range("J131").Value = 862.2317
$ws.Range("K131").Value = 1561.94109
$ws.Range("L131").Value = 2586.6951
$ws.Range("M131").Value = 3478.05891
$ws.Range("N131").Value = -12666.6951

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("N114").Value = 0
$ws.Range("H122").Value = 1678.2858
$ws.Range("J122").Value = 3980
$ws.Range("L122").Value = 11940
$ws.Range("N122").Value = -16840

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H96").Value = 16333.333
$ws.Range("J96").Value = 16333.333
$ws.Range("L96").Value = 16333.333
$ws.Range("N96").Value = -21825.333
$ws.Range("H100").Value = 2278.75
$ws.Range("I100").Value = 2133.3333
$ws.Range("J100").Value = 2366
$ws.Range("K100").Value = 2133.3333
$ws.Range("L100").Value = 2366
$ws.Range("M100").Value = -1592.3333
$ws.Range("N100").Value = -3448
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("N108").Value = 0
$ws.Range("H110").Value = 24000
$ws.Range("J110").Value = 24000
$ws.Range("L110").Value = 24000
$ws.Range("N110").Value = -32180
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("N114").Value = 0
$ws.Range("H132").Value = 4032.08
$ws.Range("I132").Value = 5192.615
$ws.Range("J132").Value = 2774.8333
$ws.Range("K132").Value = 15577.845
$ws.Range("L132").Value = 8324.499899999999
$ws.Range("M132").Value = -13047.845
$ws.Range("N132").Value = -13384.4999
$ws.Range("H136").Value = 2956.8572
$ws.Range("I136").Value = 1993.4615
$ws.Range("J136").Value = 4522.375
$ws.Range("K136").Value = 5980.3845
$ws.Range("L136").Value = 13567.125
$ws.Range("M136").Value = -3430.3845
$ws.Range("N136").Value = -18667.125

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H114").Value = 30000
$ws.Range("J114").Value = 30000
$ws.Range("L114").Value = 30000
$ws.Range("N114").Value = -38678
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("N119").Value = 0
$ws.Range("H132").Value = 3492.9531
$ws.Range("I132").Value = 1823.6383
$ws.Range("J132").Value = 8108.1177
$ws.Range("K132").Value = 5470.9149
$ws.Range("L132").Value = 24324.3531
$ws.Range("M132").Value = -2940.9149
$ws.Range("N132").Value = -29384.3531
$ws.Range("H136").Value = 23940.84
$ws.Range("I136").Value = 67446.13
$ws.Range("J136").Value = 5295.7144
$ws.Range("K136").Value = 202338.39
$ws.Range("L136").Value = 15887.1432
$ws.Range("M136").Value = -199788.39
$ws.Range("N136").Value = -20987.1432
